$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing Test RMSE (column E) values for the first DataSetup (rows 2-5) ---
$ws.Range("E2").Value = 0.1061506
$ws.Range("E3").Value = 0.1116287
$ws.Range("E4").Value = 0.14699709999999999
$ws.Range("E5").Value = 0.1053959

# --- Append new rows (6-9) for the second DataSetup ---
$ws.Range("A6").Value = "lm"
$ws.Range("B6").Value = 2
$ws.Range("D6").Value = 0.093469170000000004
$ws.Range("E6").Value = 0.12828000000000001

$ws.Range("A7").Value = "gbm"
$ws.Range("B7").Value = 2
$ws.Range("D7").Value = 0.093826179999999995
$ws.Range("E7").Value = 0.1142764

$ws.Range("A8").Value = "rf"
$ws.Range("B8").Value = 2
$ws.Range("D8").Value = 0.053008699999999999
$ws.Range("E8").Value = 0.13004679999999999

$ws.Range("A9").Value = "MARS"
$ws.Range("B9").Value = 2
$ws.Range("D9").Value = 0.1048269
$ws.Range("E9").Value = 0.1145216

# Match the active selection left after this session's last edit
$ws.Range("D9").Select()
